# Replace the empty "Content Placeholder 2" shape on slide 9 with a 2x6
# table ("Content Placeholder 4") listing tasks and implementers.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)

# Locate the empty content placeholder (idx=1) that is being replaced.
$ph = $s.Shapes.Item(2)
$phLeft = $ph.Left
$phTop = $ph.Top
$phWidth = $ph.Width
$phHeight = $ph.Height

# EMU -> point helper (PowerPoint COM geometry is expressed in points).
$emuPerPt = 12700

# Create the replacement table over the placeholder's location/size.
$tbl = $s.Shapes.AddTable(6, 2, $phLeft, $phTop, $phWidth, $phHeight)
$tbl.Name = "Content Placeholder 4"

# Remove the original empty placeholder shape now that the table stands in.
$ph.Delete()

# The table was appended at the end of the z-order; move it back one step
# so it resumes the original placeholder's slot (Title, Table, SlideNumber).
$tbl.ZOrder(3)

# Position/size the table to match the authored layout.
$tbl.Left = 1154083 / $emuPerPt
$tbl.Top = 2209800 / $emuPerPt
$tbl.Width = 10058400 / $emuPerPt

$tb = $tbl.Table

$rows = @(
    @("TASK", "IMPLEMENTER"),
    @("1.  Display Objects", "Tạ Đức Duy"),
    @("2.  Moving ", "Tạ Đức Duy"),
    @("3.  Inventory", "Nguyễn Công Duy"),
    @("4.  Game Combat", "Trương Thanh Hùng"),
    @("5.  ………...............", "……………………………….")
)

for ($r = 1; $r -le 6; $r++) {
    for ($c = 1; $c -le 2; $c++) {
        $cell = $tb.Rows.Item($r).Cells($c).Shape.TextFrame.TextRange
        $cell.Text = $rows[$r - 1][$c - 1]
        $cell.Font.Size = 28
        if ($r -eq 1) {
            $cell.ParagraphFormat.Alignment = 2
            $cell.Font.Color.ObjectThemeColor = 1
        }
    }
}

# Match the authored row heights (text at 28pt renders taller than this
# nominal height, but the saved row markup keeps the smaller value).
for ($r = 1; $r -le 6; $r++) {
    $tb.Rows.Item($r).Height = 370840 / $emuPerPt
}
